$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cols = @(3, 4, 5, 6, 7, 9, 12, 13, 14)
$data = @(
    @(5.911987446456305, 11.39208510656983, 8.76604078345289, 57.72565733924125, 3.755683001120029, 46.60782760093369, 8.007255930706178, 36.13320484868541, 16.98307495717732),
    @(5.671382409986819, 11.36985823385049, 8.666929438712158, 56.93511831928879, 3.765073379102368, 46.18505825166956, 8.015388637497136, 35.27062507874101, 16.86861584555926),
    @(5.517687811008257, 11.36091050373927, 8.605210685346936, 56.471822062941, 3.771104149998658, 45.94459211584066, 8.020960620617643, 34.73886209376798, 16.80050501527966),
    @(5.453600773689383, 11.35843347574352, 8.579854244623565, 56.2887479125942, 3.77362889143155, 45.8514480767581, 8.023377018532434, 34.5219250703665, 16.77331951280883),
    @(5.442872689386066, 11.35809239622755, 8.57563175696805, 56.25869876431394, 3.774052193440983, 45.83627534517071, 8.023787072657672, 34.48589728821162, 16.76884054273886),
    @(5.516829340201617, 11.36087238088512, 8.60486953622326, 56.46932968192272, 3.771137926974677, 45.94331626839198, 8.020992618378603, 34.735936973183, 16.80013604053955),
    @(5.830288691768554, 11.3834387118509, 8.732053483420552, 57.44859759114876, 3.758866101753325, 46.45811733660151, 8.009940213373682, 35.83640581397105, 16.94317487982781),
    @(6.396059729629362, 11.46548131503017, 8.974103717460993, 59.53625569461074, 3.736880718529984, 47.61697878009698, 7.992842811848853, 37.96500171661049, 17.23973913248815),
    @(6.780147294445873, 11.54939348989755, 9.14680268116645, 61.16053136378863, 3.721962763485435, 48.5556157812849, 7.983053025115115, 39.49525087564956, 17.46586996553227),
    @(6.947748225615999, 11.59279928626071, 9.224120770289465, 61.91649397761852, 3.715437003736395, 49.00056959116281, 7.979197309628323, 40.1811013516983, 17.57019663390583),
    @(7.010172766985795, 11.60999446280219, 9.253209538016552, 62.20498841023464, 3.713002732781433, 49.1715504741078, 7.977822865460185, 40.43911242167114, 17.6098836085485),
    @(6.996775130422186, 11.60625736925742, 9.246953370197446, 62.14276049937316, 3.713525364915007, 49.13461768479871, 7.978115073230311, 40.38362443218932, 17.60132873232378),
    @(6.952904909790962, 11.59419871734363, 9.226517788952924, 61.94018523874541, 3.715235998703368, 49.01458702919322, 7.979082518938927, 40.20236350807328, 17.5734582492869),
    @(6.925896921771622, 11.58691135093864, 9.213975350622432, 61.81638543501438, 3.716288598800774, 48.94138565098317, 7.979686249606043, 40.09110734906683, 17.55640944891586),
    @(6.769049225163744, 11.54666265654571, 9.141723781364821, 61.1114515407096, 3.722394430903116, 48.52689033974802, 7.983317005517647, 39.45020231615658, 17.45907884317127),
    @(6.670989935895243, 11.52331638382888, 9.09707326307614, 60.68320838529145, 3.726206480579688, 48.27714478079798, 7.985697219143288, 39.05422920447746, 17.3997226168695),
    @(6.613919257178713, 11.51038059262258, 9.071275057871235, 60.43851478506355, 3.72842362408306, 48.13519469472594, 7.987122533469757, 38.82552206185348, 17.36572151819152),
    @(6.594481720180075, 11.50608512344093, 9.062520580340456, 60.35595098090688, 3.729178544892463, 48.08742713841973, 7.987614795889887, 38.74792920622686, 17.35423409331792),
    @(6.681497955438687, 11.52575063039879, 9.101838516424356, 60.72862952617925, 3.725798144268106, 48.30355567577605, 7.985438019028745, 39.0964816597421, 17.40602700931673),
    @(6.96581907055482, 11.59772001543234, 9.232525455161319, 61.99962784844334, 3.714732547404036, 49.04977625378495, 7.978796035219307, 40.25565226032531, 17.5816398119703),
    @(7.145558198470958, 11.64917809686824, 9.31682510850422, 62.84317963968283, 3.707715307494948, 49.55191978595474, 7.974954090559898, 41.00320231684741, 17.69745293450005),
    @(7.050189511979153, 11.62130779305694, 9.27193808567435, 62.39185632452382, 3.71144107935823, 49.28262776820157, 7.976959058679564, 40.60520992477829, 17.63555570828705),
    @(6.676749443431048, 11.52464859328398, 9.099684543316371, 60.70808991065589, 3.72598267350392, 48.29161021969197, 7.98555502615263, 39.07738259921561, 17.40317640634704),
    @(6.248442736156909, 11.43917389342796, 8.909478471434321, 58.95473055594679, 3.742609095152095, 47.287845395247, 7.996980043431106, 37.39399076326109, 17.15794906581295)
)
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $rowNum = 2 + $i
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($rowNum, $cols[$j]).Value = $row[$j]
    }
}
